$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds what used to be row 3's data (Пирсинг / 63266965)
$ws.Range("A2").Value = 63266965
$ws.Range("B2").Value = "Пирсинг"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# New column E: header + INN value for row 2
$ws.Range("E1").Value = "inn"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "312333011857"

# Remove old row 3 data (the "Пирсинг обманка" row), leaving just the F3 style cell
$ws.Range("A3:D3").Clear()

# Column E width (target OOXML width 55.54296875 chars; engine quantizes
# ColumnWidth to 1/6-char pixel steps, so feed the pre-image that lands
# on the nearest representable step)
$ws.Columns.Item(5).ColumnWidth = 54.709635416666664

# Selection
$ws.Range("E6").Select()
